$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Hello"
$ws.Range("B6").Value = "World"
$ws.Range("C6").Value = "2025-10-01T18:31:32.283Z"
